$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the Role column (B) for the new rows first
$ws.Range("B3").Value = "DEPARTMENT_ADMIN"
$ws.Range("B4").Value = "DEPARTMENT_HEAD"
$ws.Range("B5").Value = "UPPER_MANAGEMENT"
$ws.Range("B6").Value = "MANAGEMENT"

# Row 2 - Company Admin
$ws.Range("A2").Value = "Add a new department as Company Admin"
$ws.Range("C2").Value = "Test_COMPANY_ADMIN_department"

# Row 3 - Department Admin
$ws.Range("A3").Value = "Add a new department as Department Admin"
$ws.Range("C3").Value = "Test_DEPARTMENT_ADMIN_department"

# Row 4 - Department Head
$ws.Range("A4").Value = "Add a new department as Department Head"
$ws.Range("C4").Value = "Test_DEPARTMENT_HEAD_department"

# Row 5 - Upper Management
$ws.Range("A5").Value = "Add a new department as Upper Management"
$ws.Range("C5").Value = "Test_UPPER_MANAGEMENT_department"

# Row 6 - Management
$ws.Range("A6").Value = "Add a new department as Management"
$ws.Range("C6").Value = "Test_MANAGEMENT_department"

$ws.Columns.Item(1).ColumnWidth = 32.346354166666664
$ws.Columns.Item(2).ColumnWidth = 19.166666666666668
$ws.Columns.Item(3).ColumnWidth = 28.42

$null = $ws.Range("A2").Select()
